# Training Calendar update: mark sessions through Sep 10, 2023 as Completed,
# and move the Buffer row's date from Sep 23, 2023 to Sep 17, 2023.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 10-15 (Aug 30, Aug 31, Sep 02, Sep 03, Sep 09, Sep 10) are now complete.
# Copy the existing "Completed" cell formatting (from E3) onto the Status
# cells for these rows, then update their text to "Completed".
$ws.Range("E3").Copy()
$ws.Range("E10:E15").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("E10").Value2 = "Completed"
$ws.Range("E11").Value2 = "Completed"
$ws.Range("E12").Value2 = "Completed"
$ws.Range("E13").Value2 = "Completed"
$ws.Range("E14").Value2 = "Completed"
$ws.Range("E15").Value2 = "Completed"

# Buffer row date moved up from Sep 23, 2023 to Sep 17, 2023.
$ws.Range("A17").Value2 = "Sep 17, 2023"

# Leave the cursor where the author last left it.
$ws.Range("B21").Select()
